$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("G2").Value = 0.004031
$ws.Range("M2").Value = 0.5652576666666667
$ws.Range("N2").Value = 1.695773
$ws.Range("O2").Value = 0.1543677258353495
$ws.Range("P2").Value = 0.1543677258353495
$ws.Range("Q2").Value = 0.002278553654333333
$ws.Range("R2").Value = 0.020506982889
$ws.Range("S2").Value = 0.1543677258353495
$ws.Range("T2").Value = 0.1543677258353495

# --- Row 3 updates ---
$ws.Range("G3").Value = 0.004031
$ws.Range("O3").Value = 0.175525398437655
$ws.Range("P3").Value = 0.1755253984376551
$ws.Range("Q3").Value = 0.002590852692
$ws.Range("S3").Value = 0.175525398437655
$ws.Range("T3").Value = 0.1755253984376551

# --- Row 4: replace with the former row 5 data (recomputed values) ---
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 0.004031
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.453771
$ws.Range("N4").Value = 7.361313
$ws.Range("O4").Value = 0.6701068757269955
$ws.Range("P4").Value = 0.6701068757269955
$ws.Range("Q4").Value = 0.009891150901
$ws.Range("R4").Value = 0.089020358109
$ws.Range("S4").Value = 0.6701068757269955
$ws.Range("T4").Value = 0.6701068757269955

# --- Remove the now-redundant former row 5 ---
$ws.Rows.Item(5).Delete()
